# Generate Report for Handoff
# Updates the localization-status report to reflect a fresh handoff:
#  - Status cells move from "Handed back: in sync with en-US" to "Ready for handoff"
#  - The handoff-generation timestamps are bumped forward a minute or so
#  - The now-shorter "Ready for handoff" label lets the status columns shrink

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---
$overview.Range("G2").Value = "2016-08-25 17:02:21"
$dede.Range("H2").Value     = "2016-08-25 17:02:21"
$zhcn.Range("H2").Value     = "2016-08-25 17:02:16"

# --- Column widths shrink now that the status text is shorter ---
$overview.Columns.Item(5).ColumnWidth = 16.333333
$overview.Columns.Item(6).ColumnWidth = 16.333333
$zhcn.Columns.Item(3).ColumnWidth     = 16.333333
$dede.Columns.Item(3).ColumnWidth     = 16.333333
